$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new row's data, forcing the Date column to stay text (avoid Excel's
# automatic date-serial conversion) by temporarily using a Text number format,
# then clearing formatting back to default before re-applying the same
# centered alignment used by the rest of the data rows (style index 1).
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "2025/12/19"
$ws.Range("B40").Value = "逃离鸭科夫"
$ws.Range("C40").Value = 1352

$ws.Range("A40:C40").ClearFormats()
$ws.Range("A40:C40").HorizontalAlignment = -4108
$ws.Range("A40:C40").VerticalAlignment = -4108
